$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the six "Future Work" list paragraphs (from "Future Work" up to,
#    but not including, "Performance Analysis & Future Proofing...").
# ---------------------------------------------------------------------------
$startRange = $d.Content
$startRange.Find.Execute("Future Work") | Out-Null
$startPos = $startRange.Start

$endRange = $d.Content
$endRange.Find.Execute("Performance Analysis") | Out-Null
$endPos = $endRange.Start

$deleteRange = $d.Range($startPos, $endPos)
$deleteRange.Delete()

# ---------------------------------------------------------------------------
# 2. Move the hidden "_GoBack" bookmark from the end of the document to the
#    very start of the "Performance Analysis & Future Proofing..." paragraph.
#    Re-adding a bookmark with the same name removes the previous one.
# ---------------------------------------------------------------------------
$perfRange = $d.Content
$perfRange.Find.Execute("Performance Analysis") | Out-Null
$perfStart = $perfRange.Start

$bookmarkRange = $d.Range($perfStart, $perfStart)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# ---------------------------------------------------------------------------
# 3. Drop the stale <w:lastRenderedPageBreak/> cached on the run that starts
#    the "Include ODF stuff..." paragraph. Re-writing the run's text via
#    Find/Replace causes the renderer to regenerate the run without the
#    stale page-break cache marker.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Include ODF stuff", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Include ODF stuff", 2) | Out-Null
